# "Back to kW-hr for electricity units."
#
# Row 4 (the LED "Lamp" upgrade case) had its electricity-related figures
# temporarily expressed in W-hr; this reverts the affected cells back to
# kW-hr-based units/values (a factor-of-1000 rescale), restores the
# "kW-hr" label in F4, and mirrors the stray formatted-but-empty cell that
# shows up in J5 after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F4: service unit label goes back to "kW-hr" (was "W-hr").
$ws.Range("F4").Value = "kW-hr"

# Numeric cells that move from W-hr-based to kW-hr-based values (x1000).
$ws.Range("G4").Value = 3.6
$ws.Range("J4").Value = 0.1355
$ws.Range("N4").Value = 8833.3333333333303
$ws.Range("O4").Value = 81800

# Mirror the extra (empty, formatted-like-J4) cell that appears at J5.
$ws.Range("J4").Copy()
$ws.Range("J5").PasteSpecial(-4122)

# Selection/view moved to O4 (scrolled so column N is first visible).
$ws.Range("O4").Select()
